$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing "Jogador Que Escolhe Aleatório" values for rows 3 and 5
$ws.Range("F3").Value = 85.57
$ws.Range("F5").Value = 102.04

# Update the selected cell / active range to match the new selection (F4)
$ws.Range("F4").Select()
